# Updates cryptos.xlsx cell values to the new coinranking.com snapshot
# (commit: "Updated symbol list on Sun Jan 29 16:45:00 UTC 2023 with GitHub Actions").
# Every target cell is plain text (prices/percentages stored as strings, not numbers),
# so we briefly force a text number-format before writing the value, then restore the
# cell to the default "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"   # force text interpretation (no numeric/percent auto-parse)
    $cell.Value = $val
    $cell.Style = "Normal"     # drop the temporary text format again
}

Set-TextValue $ws 'D2' '317.41'
Set-TextValue $ws 'E2' '3.31%'
Set-TextValue $ws 'D3' '39.64'
Set-TextValue $ws 'E3' '1.63%'
Set-TextValue $ws 'D4' '5.148'
Set-TextValue $ws 'E4' '0.88%'
Set-TextValue $ws 'D5' '0.08216'
Set-TextValue $ws 'E5' '1.95%'
Set-TextValue $ws 'D6' '2.004'
Set-TextValue $ws 'E6' '4.10%'
Set-TextValue $ws 'D7' '8.278'
Set-TextValue $ws 'E7' '4.24%'
Set-TextValue $ws 'B8' 'MXToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '0.9329'
Set-TextValue $ws 'E8' '0.33%'
Set-TextValue $ws 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D9' '0.1413'
Set-TextValue $ws 'E9' '-2.96%'
Set-TextValue $ws 'B10' 'WazirX'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1979'
Set-TextValue $ws 'E10' '2.42%'
Set-TextValue $ws 'B11' 'MandalaExchangeToken'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.09085'
Set-TextValue $ws 'E11' '0.85%'
Set-TextValue $ws 'B12' 'BitrueCoin'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D12' '0.03545'
Set-TextValue $ws 'E12' '1.52%'
Set-TextValue $ws 'B13' 'BitMartToken'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D13' '0.09796'
Set-TextValue $ws 'E13' '0.11%'
Set-TextValue $ws 'B14' 'BitForexToken'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D14' '0.001401'
Set-TextValue $ws 'E14' '0.10%'
Set-TextValue $ws 'B15' 'TigerCash'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D15' '0.006411'
Set-TextValue $ws 'E15' '5.52%'
Set-TextValue $ws 'B16' 'LEO'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D16' '3.679'
Set-TextValue $ws 'E16' '-1.97%'
Set-TextValue $ws 'B17' 'GateToken'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D17' '4.286'
Set-TextValue $ws 'E17' '2.33%'
Set-TextValue $ws 'D18' '3.205'
Set-TextValue $ws 'E18' '-6.16%'
Set-TextValue $ws 'D19' '0.3490'
Set-TextValue $ws 'E19' '1.29%'
Set-TextValue $ws 'D20' '0.1295'
Set-TextValue $ws 'E20' '-0.44%'
Set-TextValue $ws 'D21' '4.904'
Set-TextValue $ws 'E21' '2.57%'
Set-TextValue $ws 'D22' '0.2454'
Set-TextValue $ws 'E22' '-1.95%'
Set-TextValue $ws 'D23' '0.04332'
Set-TextValue $ws 'E23' '-0.95%'
Set-TextValue $ws 'D24' '0.001226'
Set-TextValue $ws 'E24' '-0.79%'
Set-TextValue $ws 'D25' '0.004762'
Set-TextValue $ws 'E25' '11.22%'
Set-TextValue $ws 'D26' '0.0001301'
Set-TextValue $ws 'E26' '0.10%'
Set-TextValue $ws 'D27' '0.0004002'
Set-TextValue $ws 'E27' '-10.03%'
Set-TextValue $ws 'D39' '0.02216'
Set-TextValue $ws 'E39' '7.51%'
Set-TextValue $ws 'D40' '0.05220'
Set-TextValue $ws 'E40' '3.40%'
Set-TextValue $ws 'B41' 'KickToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D41' '0.007536'
Set-TextValue $ws 'E41' '1.18%'
Set-TextValue $ws 'B42' 'Dexo'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
Set-TextValue $ws 'D42' '0.01028'
Set-TextValue $ws 'E42' '1.83%'
Set-TextValue $ws 'D43' '0.1376'
Set-TextValue $ws 'E43' '1.82%'
Set-TextValue $ws 'E44' '0.57%'
Set-TextValue $ws 'D45' '0.009852'
Set-TextValue $ws 'E45' '8.72%'
Set-TextValue $ws 'D46' '0.00006640'
Set-TextValue $ws 'E46' '7.50%'
Set-TextValue $ws 'E47' '0.08%'
Set-TextValue $ws 'B48' 'BOLO'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws 'D48' '0.002766'
Set-TextValue $ws 'E48' '-1.09%'
Set-TextValue $ws 'B49' 'CoinbaseStockToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws 'D49' '0.001200'
Set-TextValue $ws 'E49' '-24.94%'
Set-TextValue $ws 'E50' '0.08%'
Set-TextValue $ws 'E51' '0.08%'

Write-Output "Applied $(102) cell updates"
